# Replace the shared-string "guia" tracking numbers (A2:A5, padded with trailing
# spaces) with a cleaned-up, longer list of numeric guide/tracking numbers (A2:A70).
# Values are written as real numbers (not text) and right-aligned / vertically
# centered in black Calibri, matching the "limpieza del dato a buscar" cleanup.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    231778496,
    2161788829,
    2162143803,
    2165877749,
    2167602235,
    2169277038,
    9137160126,
    9147007918,
    522760109,
    2095307865,
    2111368424,
    2111368431,
    2123282428,
    2167602461,
    2169146641,
    2169933592,
    2170422314,
    9154439920,
    9155305065,
    2089161788,
    2089161812,
    2089161922,
    2089161941,
    2089166520,
    2092670202,
    2104592813,
    2107373906,
    2107376639,
    2111374129,
    2111381675,
    2111384338,
    2111408816,
    2124902530,
    2131458302,
    2131458587,
    2132213674,
    2132467190,
    2137035982,
    2143590008,
    2145353943,
    2154087258,
    2154087259,
    2159691695,
    2160511138,
    2161855791,
    2170295802,
    2170295807,
    2170339631,
    2170506586,
    2170672262,
    2173706594,
    2174746520,
    3006323249,
    3014462180,
    3041651686,
    3041652089,
    3041989991,
    3041991060,
    3044095813,
    3044139383,
    3044139471,
    3044415861,
    3044418127,
    9147225713,
    9153467569,
    9154939023,
    9155305185,
    9157577296,
    9157957057
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Format the first data cell, then fan the format out to the rest of the column via
# copy/paste-special so every cell lands on the same single style record instead of
# each cell accreting its own intermediate (font-only / partially-aligned) style.
$firstCell = $ws.Cells.Item(2, 1)
$firstCell.Font.Color = 0
$firstCell.HorizontalAlignment = -4152   # xlRight
$firstCell.VerticalAlignment = -4108     # xlCenter
$firstCell.Copy()
for ($row = 3; $row -le 70; $row++) {
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

$ws.Range("A2").Select()
